# Update the "想去人数" (want-to-go count) column F values across the four
# worksheets to match the freshly re-scraped figures.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)  # 展览
$ws2 = $wb.Worksheets.Item(2)  # 演出
$ws3 = $wb.Worksheets.Item(3)  # 本地生活
$ws4 = $wb.Worksheets.Item(4)  # 全部类型

# 展览 (sheet 1)
$ws1.Range("F2").Value = 21
$ws1.Range("F7").Value = 3952
$ws1.Range("F9").Value = 801
$ws1.Range("F10").Value = 2405
$ws1.Range("F11").Value = 376
$ws1.Range("F12").Value = 57
$ws1.Range("F13").Value = 243
$ws1.Range("F17").Value = 4025
$ws1.Range("F19").Value = 240
$ws1.Range("F22").Value = 253
$ws1.Range("F24").Value = 290

# 演出 (sheet 2)
$ws2.Range("F7").Value = 137
$ws2.Range("F23").Value = 89

# 本地生活 (sheet 3)
$ws3.Range("F2").Value = 6407
$ws3.Range("F4").Value = 2143
$ws3.Range("F6").Value = 30

# 全部类型 (sheet 4)
$ws4.Range("F3").Value = 2143
$ws4.Range("F6").Value = 21
$ws4.Range("F14").Value = 30
$ws4.Range("F16").Value = 3952
$ws4.Range("F17").Value = 137
$ws4.Range("F21").Value = 801
$ws4.Range("F22").Value = 2405
$ws4.Range("F23").Value = 376
$ws4.Range("F24").Value = 57
$ws4.Range("F26").Value = 243
$ws4.Range("F36").Value = 240
$ws4.Range("F38").Value = 361
$ws4.Range("F39").Value = 253
$ws4.Range("F48").Value = 89
$ws4.Range("F49").Value = 290

$wb.Save()
